# Insert a new row at 520, which shifts rows 520:593 down to 521:594,
# preserving all of their existing data/formatting (matches the diff's
# "shift" pattern where every row 520-592 takes on the values previously
# held by the next row, and old row 593 becomes new row 594).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("520:520").Insert()

# Populate the newly inserted row 520 with the new data record.
$ws.Range("A520").Value = 8
$ws.Range("B520").Value = "Terminal La Palmera de La Serena"
$ws.Range("C520").Value = "Coquimbo"
$ws.Range("D520").Value = 44984
$ws.Range("E520").Value = 4
$ws.Range("F520").Value = 100114001
$ws.Range("G520").Value = "Papa"
$ws.Range("H520").Value = "Asterix"
$ws.Range("I520").Value = "1a (cosecha)"
$ws.Range("J520").Value = 1800
$ws.Range("K520").Value = 11500
$ws.Range("L520").Value = 12000
$ws.Range("M520").Value = 11750
$ws.Range("N520").Value = "$/saco 25 kilos"
$ws.Range("O520").Value = "Región del Maule"
$ws.Range("P520").Value = 470
$ws.Range("Q520").Value = 25
$ws.Range("R520").Value = "Hortaliza"
